$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells are stored as text so values like
# "1.003" are not auto-converted into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.062.72"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.888.96"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.20"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5038"
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3895"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09174"
$ws.Range("E9").Value = "  -6.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.125"
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.80"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.373"
$ws.Range("E12").Value = "  -2.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.78"
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.900.97"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.270"
$ws.Range("E15").Value = "  -4.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.36"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001105"
$ws.Range("E18").Value = "  -3.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06666"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.80"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.193"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.128.40"
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.38"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.324"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.116.80"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.536"
$ws.Range("E27").Value = "  -7.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.15"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.79"
$ws.Range("E29").Value = "  -2.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.85"
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.072"
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.591"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.612"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.484"
$ws.Range("E35").Value = "  -3.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.346"
$ws.Range("E36").Value = "  +13.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06591"
$ws.Range("E37").Value = "  -3.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02396"
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2201"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.214"
$ws.Range("E40").Value = "  -4.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6431"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("E42").Value = "  -4.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.951"
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.33"
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6047"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.304"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.685"
$ws.Range("E48").Value = "  -3.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.994"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.196"
$ws.Range("E51").Value = "  -1.55%  "
